$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns (STRASSE, HAUSNR) right before the existing PLZ column (J)
$ws.Range("J1:K1").EntireColumn.Insert()

# Fill in the new columns - order chosen to match the shared-string table ordering
$ws.Range("J2").Value2 = "Ackerstrasse"
$ws.Range("K1").Value2 = "HAUSNR"
$ws.Range("J1").Value2 = "STRASSE"
$ws.Range("J3").Value2 = "Denzingsteig"
$ws.Range("J4").Value2 = "Eichenweg"
$ws.Range("K2").Value2 = 11
$ws.Range("K3").Value2 = 22
$ws.Range("K4").Value2 = 33

# Match the column width used for the neighbouring STRASSE/HAUSNR columns
$ws.Range("J1:K1").ColumnWidth = 22.33

# Update the active selection to the newly added row-3 cells
$ws.Range("J3:K3").Select()
